# Update "想去人数" (want-to-go count) values in column F
# for the "展览" and "全部类型" sheets, as produced by the
# gh-pages data regeneration at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value  = 8447
    $ws.Range("F3").Value  = 8086
    $ws.Range("F9").Value  = 140
    $ws.Range("F12").Value = 730
    $ws.Range("F13").Value = 180
    $ws.Range("F14").Value = 2193
    $ws.Range("F19").Value = 141
}
